$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.767168770400019
$ws.Range("C2").Value = 0.3573878827427563
$ws.Range("D2").Value = 0.009189471803228599
$ws.Range("E2").Value = 0.04535142159905181
$ws.Range("F2").Value = 3.696417607250041
$ws.Range("I2").Value = 2.157564468342706
$ws.Range("J2").Value = 0.1037402433238146
$ws.Range("L2").Value = 0.3924357081746805
$ws.Range("M2").Value = 0.5907380481850169
$ws.Range("N2").Value = 2.567041630754666
$ws.Range("B3").Value = 2.678631307482931
$ws.Range("C3").Value = 0.3305849875997637
$ws.Range("D3").Value = 0.008639454830511539
$ws.Range("E3").Value = 0.044987651525072
$ws.Range("F3").Value = 3.688666605529249
$ws.Range("I3").Value = 2.160225284905678
$ws.Range("J3").Value = 0.1037164097845853
$ws.Range("L3").Value = 0.3902464738432059
$ws.Range("M3").Value = 0.5774678069625523
$ws.Range("N3").Value = 2.590871131174119
$ws.Range("B4").Value = 2.625770047447475
$ws.Range("C4").Value = 0.3143289265098588
$ws.Range("D4").Value = 0.008297585276089592
$ws.Range("E4").Value = 0.04475934156067929
$ws.Range("F4").Value = 3.685739604632062
$ws.Range("I4").Value = 2.162872607125458
$ws.Range("J4").Value = 0.1037040203126169
$ws.Range("L4").Value = 0.3890673088530008
$ws.Range("M4").Value = 0.5696262734234452
$ws.Range("N4").Value = 2.606282021093925
$ws.Range("B5").Value = 2.604605926700799
$ws.Range("C5").Value = 0.3077545846122689
$ws.Range("D5").Value = 0.008157181559234061
$ws.Range("E5").Value = 0.04466504598543164
$ws.Range("F5").Value = 3.685006833174327
$ws.Range("I5").Value = 2.164205901958354
$ws.Range("J5").Value = 0.1036995374742373
$ws.Range("L5").Value = 0.3886283279267033
$ws.Range("M5").Value = 0.5665078430751151
$ws.Range("N5").Value = 2.612757886682473
$ws.Range("B6").Value = 2.601114421895716
$ws.Range("C6").Value = 0.3066659361981863
$ws.Range("D6").Value = 0.008133800573290273
$ws.Range("E6").Value = 0.04464931198927236
$ws.Range("F6").Value = 3.684912919966493
$ws.Range("I6").Value = 2.164442654093058
$ws.Range("J6").Value = 0.1036988273336164
$ws.Range("L6").Value = 0.3885579454863688
$ws.Range("M6").Value = 0.5659946843121944
$ws.Range("N6").Value = 2.613845019354045
$ws.Range("B7").Value = 2.62548309340599
$ws.Range("C7").Value = 0.3142400601812767
$ws.Range("D7").Value = 0.008295696209277281
$ws.Range("E7").Value = 0.04475807496404727
$ws.Range("F7").Value = 3.685727860625377
$ws.Range("I7").Value = 2.162889558505675
$ws.Range("J7").Value = 0.1037039575612777
$ws.Range("L7").Value = 0.3890612203617465
$ws.Range("M7").Value = 0.5695839051754348
$ws.Range("N7").Value = 2.606368564371792
$ws.Range("B8").Value = 2.736329342627357
$ws.Range("C8").Value = 0.3481042622232735
$ws.Range("D8").Value = 0.009000662188476838
$ws.Range("E8").Value = 0.04522701464574919
$ws.Range("F8").Value = 3.693364361819263
$ws.Range("I8").Value = 2.158271316009156
$ws.Range("J8").Value = 0.1037315606174278
$ws.Range("L8").Value = 0.3916466290038798
$ws.Range("M8").Value = 0.586098858369219
$ws.Range("N8").Value = 2.575096140284032
$ws.Range("B9").Value = 2.965633683675833
$ws.Range("C9").Value = 0.4161294350551543
$ws.Range("D9").Value = 0.01035213267090285
$ws.Range("E9").Value = 0.0461078570859712
$ws.Range("F9").Value = 3.722914251355007
$ws.Range("I9").Value = 2.157276658493075
$ws.Range("J9").Value = 0.1038034384508837
$ws.Range("L9").Value = 0.3980251022053238
$ws.Range("M9").Value = 0.6209188409758113
$ws.Range("N9").Value = 2.519964152420577
$ws.Range("B10").Value = 3.141436773030534
$ws.Range("C10").Value = 0.4671303913287375
$ws.Range("D10").Value = 0.01132903325937562
$ws.Range("E10").Value = 0.04673220311902693
$ws.Range("F10").Value = 3.753569981616181
$ws.Range("I10").Value = 2.161490814973732
$ws.Range("J10").Value = 0.1038670009276093
$ws.Range("L10").Value = 0.403508869842014
$ws.Range("M10").Value = 0.647992586528261
$ws.Range("N10").Value = 2.483244440010111
$ws.Range("B11").Value = 3.223020637216678
$ws.Range("C11").Value = 0.4905624910067559
$ws.Range("D11").Value = 0.01177059343518394
$ws.Range("E11").Value = 0.04701145940366036
$ws.Range("F11").Value = 3.769472287082749
$ws.Range("I11").Value = 2.16448858152286
$ws.Range("J11").Value = 0.1038982398577843
$ws.Range("L11").Value = 0.406176810923867
$ws.Range("M11").Value = 0.6606349379672309
$ws.Range("N11").Value = 2.467364006691142
$ws.Range("B12").Value = 3.254146587533853
$ws.Range("C12").Value = 0.4994694771293666
$ws.Range("D12").Value = 0.01193744266516461
$ws.Range("E12").Value = 0.04711653529518767
$ws.Range("F12").Value = 3.775776478165653
$ws.Range("I12").Value = 2.165779689439148
$ws.Range("J12").Value = 0.1039104022071324
$ws.Range("L12").Value = 0.4072120051474712
$ws.Range("M12").Value = 0.6654692879338455
$ws.Range("N12").Value = 2.461469172150906
$ws.Range("B13").Value = 3.247432736640235
$ws.Range("C13").Value = 0.4975496893531499
$ws.Range("D13").Value = 0.01190152394362798
$ws.Range("E13").Value = 0.04709393502967973
$ws.Range("F13").Value = 3.77440618581241
$ws.Range("I13").Value = 2.165494683460793
$ws.Range("J13").Value = 0.1039077680451985
$ws.Range("L13").Value = 0.4069879503575464
$ws.Range("M13").Value = 0.6644260347499866
$ws.Range("N13").Value = 2.462733444097495
$ws.Range("B14").Value = 3.225576737791869
$ws.Range("C14").Value = 0.4912945948389051
$ws.Range("D14").Value = 0.01178432716183764
$ws.Range("E14").Value = 0.0470201174653484
$ws.Range("F14").Value = 3.769985272855251
$ws.Range("I14").Value = 2.164591673938517
$ws.Range("J14").Value = 0.1038992337983942
$ws.Range("L14").Value = 0.4062614780188483
$ws.Range("M14").Value = 0.6610317213680119
$ws.Range("N14").Value = 2.466876653528267
$ws.Range("B15").Value = 3.212219520704764
$ws.Range("C15").Value = 0.4874675792034395
$ws.Range("D15").Value = 0.01171249529203067
$ws.Range("E15").Value = 0.04697481487716715
$ws.Range("F15").Value = 3.767314132318603
$ws.Range("I15").Value = 2.164058875112261
$ws.Range("J15").Value = 0.1038940496326859
$ws.Range("L15").Value = 0.4058197351275936
$ws.Range("M15").Value = 0.6589587238595769
$ws.Range("N15").Value = 2.469429961998671
$ws.Range("B16").Value = 3.136137507340095
$ws.Range("C16").Value = 0.4656037468789691
$ws.Range("D16").Value = 0.01130012326481378
$ws.Range("E16").Value = 0.04671385829664043
$ws.Range("F16").Value = 3.752570185415493
$ws.Range("I16").Value = 2.161316694366377
$ws.Range("J16").Value = 0.1038650060398716
$ws.Range("L16").Value = 0.4033379996944859
$ws.Range("M16").Value = 0.6471729493325071
$ws.Range("N16").Value = 2.484298862387369
$ws.Range("B17").Value = 3.089876312570141
$ws.Range("C17").Value = 0.4522506194144285
$ws.Range("D17").Value = 0.01104645375632174
$ws.Range("E17").Value = 0.04655255901742716
$ws.Range("F17").Value = 3.744027112404581
$ws.Range("I17").Value = 2.159911616093154
$ws.Range("J17").Value = 0.1038477830040598
$ws.Range("L17").Value = 0.4018599159927447
$ws.Range("M17").Value = 0.6400263721896025
$ws.Range("N17").Value = 2.493631603885426
$ws.Range("B18").Value = 3.063419628466704
$ws.Range("C18").Value = 0.4445920159603816
$ws.Range("D18").Value = 0.01090027985980058
$ws.Range("E18").Value = 0.0464593349857827
$ws.Range("F18").Value = 3.739297494425244
$ws.Range("I18").Value = 2.159205144083103
$ws.Range("J18").Value = 0.103838095597375
$ws.Range("L18").Value = 0.4010260789923024
$ws.Range("M18").Value = 0.6359465670525992
$ws.Range("N18").Value = 2.499077046893511
$ws.Range("B19").Value = 3.054487864600389
$ws.Range("C19").Value = 0.4420026679844113
$ws.Range("D19").Value = 0.01085074008192422
$ws.Range("E19").Value = 0.04642769344803455
$ws.Range("F19").Value = 3.737727723507987
$ws.Range("I19").Value = 2.158983395022261
$ws.Range("J19").Value = 0.1038348532337778
$ws.Range("L19").Value = 0.4007465594536939
$ws.Range("M19").Value = 0.634570489840975
$ws.Range("N19").Value = 2.500934080987648
$ws.Range("B20").Value = 3.094785210208158
$ws.Range("C20").Value = 0.4536698269426438
$ws.Range("D20").Value = 0.01107348489642135
$ws.Range("E20").Value = 0.04656977597607259
$ws.Range("F20").Value = 3.744917473505666
$ws.Range("I20").Value = 2.160050660718383
$ws.Range("J20").Value = 0.1038495937877517
$ws.Range("L20").Value = 0.4020155718946228
$ws.Range("M20").Value = 0.6407839582045014
$ws.Range("N20").Value = 2.492630094762262
$ws.Range("B21").Value = 3.231990077322337
$ws.Range("C21").Value = 0.4931309477615287
$ws.Range("D21").Value = 0.01181876008401161
$ws.Range("E21").Value = 0.04704181762337001
$ws.Range("F21").Value = 3.771276132134574
$ws.Range("I21").Value = 2.164852674005658
$ws.Range("J21").Value = 0.1039017314904758
$ws.Range("L21").Value = 0.4064741848837627
$ws.Range("M21").Value = 0.6620274385487406
$ws.Range("N21").Value = 2.465656467338313
$ws.Range("B22").Value = 3.323013341181593
$ws.Range("C22").Value = 0.5191180643458324
$ws.Range("D22").Value = 0.01230376454229898
$ws.Range("E22").Value = 0.04734641064543244
$ws.Range("F22").Value = 3.790149142401077
$ws.Range("I22").Value = 2.168900125759279
$ws.Range("J22").Value = 0.1039377459707715
$ws.Range("L22").Value = 0.4095332890088486
$ws.Range("M22").Value = 0.6761850534748461
$ws.Range("N22").Value = 2.448719907082314
$ws.Range("B23").Value = 3.274308677796739
$ws.Range("C23").Value = 0.505230069745437
$ws.Range("D23").Value = 0.01204508254370396
$ws.Range("E23").Value = 0.04718419759761439
$ws.Range("F23").Value = 3.77992532776517
$ws.Range("I23").Value = 2.166656574021928
$ws.Range("J23").Value = 0.1039183473368812
$ws.Range("L23").Value = 0.4078873145606394
$ws.Range("M23").Value = 0.668603805062979
$ws.Range("N23").Value = 2.457695819566545
$ws.Range("B24").Value = 3.092565463107746
$ws.Range("C24").Value = 0.4530281462971857
$ws.Range("D24").Value = 0.01106126516805261
$ws.Range("E24").Value = 0.0465619937194317
$ws.Range("F24").Value = 3.744514374686275
$ws.Range("I24").Value = 2.159987483052845
$ws.Range("J24").Value = 0.1038487744651313
$ws.Range("L24").Value = 0.4019451502311426
$ws.Range("M24").Value = 0.6404413636828394
$ws.Range("N24").Value = 2.493082628087798
$ws.Range("B25").Value = 2.902317230790402
$ws.Range("C25").Value = 0.3975498781706506
$ws.Range("D25").Value = 0.00998954579001321
$ws.Range("E25").Value = 0.04587363032348879
$ws.Range("F25").Value = 3.713353733863883
$ws.Range("I25").Value = 2.156679580816885
$ws.Range("J25").Value = 0.1037821004900454
$ws.Range("L25").Value = 0.3961594989406763
$ws.Range("M25").Value = 0.6609188409758113
$ws.Range("N25").Value = 2.534214764090862
